# Actualizacion de las graficas generales del plan del proyecto
#
# Ciclo #4: el task #66 ("cierre del ciclo #3") se renombra/corrige a #61 en
# "task" y en su registro de tiempo en "logt"; las horas trabajadas de la
# tarea #72 ahora se calculan con SUMIF (como ya ocurria en las otras filas)
# en vez de mantenerse como entradas manuales repetidas, así que las filas
# 5 y 6 (duplicados de la tarea #72) quedan vacías.

$wb = $excel.ActiveWorkbook

$schedule = $wb.Worksheets.Item("schedule")
$task     = $wb.Worksheets.Item("task")
$logt     = $wb.Worksheets.Item("logt")

# --- logt: corrige el id de la tarea referenciada en la bitacora ---
$logt.Range("G2").Value = 61

# --- task: corrige el id y agrega formulas SUMIF que ya se usaban en
#     las demas filas (en vez de horas acumuladas manualmente) ---
$task.Range("A2").Value = 61
$task.Range("B3").Formula = "=SUMIF(logt!G:G,task!A3,logt!F:F)/60"
$task.Range("B4").Formula = "=SUMIF(logt!G:G,task!A4,logt!F:F)/60"

# Las filas 5 y 6 eran entradas duplicadas de la tarea #72 que ahora se
# consolidan en la formula de B4, asi que se limpian.
$task.Range("A5:C6").ClearContents()

# --- selecciones / hoja activa: el autor termino revisando "logt" y
#     luego "task" (que queda como la hoja activa al guardar) ---
$schedule.Range("B4:D4").Select()
$logt.Range("G3").Select()
$task.Select()
$task.Range("A3").Select()
